$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculated K (Strike#) values for rows 2-15, column G
$kValues = @{
    2  = 3
    3  = 3
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 3
    13 = 2
    14 = 0
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
